# Auto-generated Excel COM-interop script applying the scheduled-runner update
# to Sheets/Halicarnassus_Profits.xlsx (workbook with sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each affected row, set/clear the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# to the refreshed values captured by the latest price-data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 1647.5333
$ws.Range("I4").Value = 1820.6666
$ws.Range("J4").Value = 1387.8334
$ws.Range("K4").Value = 1820.6666
$ws.Range("L4").Value = 1387.8334
$ws.Range("M4").Value = -1706.6666
$ws.Range("N4").Value = -1615.8334
# Row 6
$ws.Range("H6").Value = 148
$ws.Range("J6").Value = 148
$ws.Range("L6").Value = 444
$ws.Range("N6").Value = -668
# Row 43
$ws.Range("H43").Value = 2483.6365
$ws.Range("J43").Value = 2451.2222
$ws.Range("L43").Value = 2451.2222
$ws.Range("N43").Value = -2589.2222
# Row 74
$ws.Range("H74").Value = 7810.625
$ws.Range("I74").Value = 8283.714
$ws.Range("K74").Value = 8283.714
$ws.Range("M74").Value = -7347.714
# Row 77
$ws.Range("H77").Value = 7810.625
$ws.Range("I77").Value = 8283.714
$ws.Range("K77").Value = 41418.57
$ws.Range("M77").Value = -36738.57
# Row 100
$ws.Range("H100").Value = 6548.8
$ws.Range("I100").Value = 1500
$ws.Range("J100").Value = 7811
$ws.Range("K100").Value = 1500
$ws.Range("L100").Value = 7811
$ws.Range("M100").Value = -959
$ws.Range("N100").Value = -8893
# Row 116
$ws.Range("H116").Value = 6777.4
$ws.Range("I116").Value = 6333.3335
$ws.Range("K116").Value = 6333.3335
$ws.Range("M116").Value = -2891.3335
# Row 132
$ws.Range("H132").Value = 2161.8076
$ws.Range("I132").Value = 1248.32
$ws.Range("K132").Value = 3744.96
$ws.Range("M132").Value = -1214.96
# Row 137
$ws.Range("H137").Value = 4887
$ws.Range("I137").Value = 2112.4546
$ws.Range("K137").Value = 6337.3638
$ws.Range("M137").Value = -3787.3638
# Row 138
$ws.Range("H138").Value = 3063.0952
$ws.Range("I138").Value = 1562.4286
$ws.Range("J138").Value = 3813.4285
$ws.Range("K138").Value = 4687.2858
$ws.Range("L138").Value = 11440.2855
$ws.Range("M138").Value = 452.7142000000003
$ws.Range("N138").Value = -21720.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2923
$ws.Range("I2").Value = 2923
$ws.Range("K2").Value = 2923
$ws.Range("M2").Value = -2810
# Row 32
$ws.Range("H32").Value = 17417.525
$ws.Range("I32").Value = 15086.077
$ws.Range("J32").Value = 22469
$ws.Range("K32").Value = 15086.077
$ws.Range("L32").Value = 22469
$ws.Range("M32").Value = -14799.077
$ws.Range("N32").Value = -23043
# Row 74
$ws.Range("H74").Value = 2028.6818
$ws.Range("I74").Value = 1723.3158
$ws.Range("K74").Value = 1723.3158
$ws.Range("M74").Value = -849.3158000000001
# Row 77
$ws.Range("H77").Value = 2028.6818
$ws.Range("I77").Value = 1723.3158
$ws.Range("K77").Value = 8616.579
$ws.Range("M77").Value = -4248.579
# Row 116
$ws.Range("H116").Value = 2923
$ws.Range("I116").Value = 2923
$ws.Range("K116").Value = 2923
$ws.Range("M116").Value = -629
# Row 132
$ws.Range("H132").Value = 733.3333
$ws.Range("J132").Value = 600
$ws.Range("L132").Value = 1800
$ws.Range("N132").Value = -6860

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2923
$ws.Range("I3").Value = 2923
$ws.Range("K3").Value = 2923
$ws.Range("M3").Value = -2809
# Row 86
$ws.Range("H86").Value = 6190.5557
$ws.Range("I86").Value = 5316.625
$ws.Range("K86").Value = 5316.625
$ws.Range("M86").Value = -4193.625
# Row 89
$ws.Range("H89").Value = 6190.5557
$ws.Range("I89").Value = 5316.625
$ws.Range("K89").Value = 26583.125
$ws.Range("M89").Value = -20967.125
# Row 94
$ws.Range("H94").Value = 940.6667
$ws.Range("I94").Value = 948.5
$ws.Range("J94").Value = 925
$ws.Range("K94").Value = 948.5
$ws.Range("L94").Value = 925
$ws.Range("M94").Value = -497.5
$ws.Range("N94").Value = -1827
# Row 134
$ws.Range("H134").Value = 3888.3333
$ws.Range("I134").Value = 3256.5715
$ws.Range("K134").Value = 9769.7145
$ws.Range("M134").Value = -7234.7145

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2875
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3574
# Row 31
$ws.Range("H31").Value = 5335.032
$ws.Range("I31").Value = 4709.615
$ws.Range("J31").Value = 8587.200000000001
$ws.Range("K31").Value = 4709.615
$ws.Range("L31").Value = 8587.200000000001
$ws.Range("M31").Value = -4414.615
$ws.Range("N31").Value = -9177.200000000001
# Row 34
$ws.Range("H34").Value = 5335.032
$ws.Range("I34").Value = 4709.615
$ws.Range("J34").Value = 8587.200000000001
$ws.Range("K34").Value = 4709.615
$ws.Range("L34").Value = 8587.200000000001
$ws.Range("M34").Value = -4507.615
$ws.Range("N34").Value = -8991.200000000001
# Row 58
$ws.Range("H58").Value = 3070.6667
$ws.Range("I58").Value = 1918.6316
$ws.Range("J58").Value = 7448.4
$ws.Range("K58").Value = 1918.6316
$ws.Range("L58").Value = 7448.4
$ws.Range("M58").Value = -1715.6316
$ws.Range("N58").Value = -7854.4
# Row 63
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372
# Row 66
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864
# Row 103
$ws.Range("H103").Value = 15645
$ws.Range("I103").Value = 15645
$ws.Range("K103").Value = 15645
$ws.Range("M103").Value = -14473
# Row 105
$ws.Range("H105").Value = 1648.5714
$ws.Range("I105").Value = 1136.3334
$ws.Range("K105").Value = 1136.3334
$ws.Range("M105").Value = 610.6666
# Row 107
$ws.Range("H107").Value = 399.20834
$ws.Range("I107").Value = 333.6
$ws.Range("K107").Value = 333.6
$ws.Range("M107").Value = 1586.4
# Row 113
$ws.Range("H113").Value = 2875
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
# Row 132
$ws.Range("H132").Value = 2997.889
$ws.Range("I132").Value = 2997.889
$ws.Range("K132").Value = 8993.667000000001
$ws.Range("M132").Value = -6463.667000000001
# Row 136
$ws.Range("H136").Value = 3070.6667
$ws.Range("I136").Value = 1918.6316
$ws.Range("J136").Value = 7448.4
$ws.Range("K136").Value = 5755.8948
$ws.Range("L136").Value = 22345.2
$ws.Range("M136").Value = -3205.8948
$ws.Range("N136").Value = -27445.2

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 962150.6
$ws.Range("I4").Value = 600301.4
$ws.Range("K4").Value = 1800904.2
$ws.Range("M4").Value = -1800792.2
# Row 116
$ws.Range("H116").Value = 2800
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
# Row 132
$ws.Range("H132").Value = 3212.9033
$ws.Range("I132").Value = 4200.143
$ws.Range("K132").Value = 37801.287
$ws.Range("M132").Value = -35271.287
# Row 141
$ws.Range("H141").Value = 2028.25
$ws.Range("I141").Value = 1889.4286
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 5668.2858
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -488.2857999999997
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 33245
$ws.Range("J15").Value = 30326.666
$ws.Range("L15").Value = 30326.666
$ws.Range("N15").Value = -30902.666
# Row 43
$ws.Range("H43").Value = 8000
$ws.Range("I43").Value = 8000
$ws.Range("K43").Value = 8000
$ws.Range("M43").Value = -7849
# Row 81
$ws.Range("H81").Value = 33245
$ws.Range("J81").Value = 30326.666
$ws.Range("L81").Value = 30326.666
$ws.Range("N81").Value = -32322.666
# Row 84
$ws.Range("H84").Value = 33245
$ws.Range("J84").Value = 30326.666
$ws.Range("L84").Value = 90979.99800000001
$ws.Range("N84").Value = -100963.998
# Row 126
$ws.Range("H126").Value = 3489
$ws.Range("I126").Value = 3489
$ws.Range("K126").Value = 10467
$ws.Range("M126").Value = -7997
# Row 132
$ws.Range("H132").Value = 142901.12
$ws.Range("I132").Value = 186368.5
$ws.Range("K132").Value = 559105.5
$ws.Range("M132").Value = -556575.5

$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 40
$ws.Range("H40").Value = 11666
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 11666
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 11666
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -11938
# Row 122
$ws.Range("H122").Value = 4232.6665
$ws.Range("I122").Value = 4232.6665
$ws.Range("K122").Value = 12697.9995
$ws.Range("M122").Value = -10247.9995
# Row 132
$ws.Range("H132").Value = 8291.611000000001
$ws.Range("I132").Value = 7054.3335
$ws.Range("J132").Value = 9528.888999999999
$ws.Range("K132").Value = 21163.0005
$ws.Range("L132").Value = 28586.667
$ws.Range("M132").Value = -18633.0005
$ws.Range("N132").Value = -33646.667

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3487.8928
$ws.Range("I122").Value = 2568.353
$ws.Range("K122").Value = 7705.059
$ws.Range("M122").Value = -5255.059
# Row 126
$ws.Range("H126").Value = 7000
$ws.Range("I126").Value = 4075
$ws.Range("K126").Value = 12225
$ws.Range("M126").Value = -9755
# Row 132
$ws.Range("H132").Value = 5564.7617
$ws.Range("I132").Value = 5160.7334
$ws.Range("K132").Value = 15482.2002
$ws.Range("M132").Value = -12952.2002
# Row 136
$ws.Range("H136").Value = 6135.7417
$ws.Range("J136").Value = 9139.857
$ws.Range("L136").Value = 27419.571
$ws.Range("N136").Value = -32519.571

Write-Host "Applied Halicarnassus_Profits price refresh."
